# Daily attendance processing - 2025-12-29 15:02:07
# Reorders the "Recorded By" (column G) contributor lists on the
# "Session Analysis Results" sheet for specific rows: the first
# name/email in the comma-separated list is moved to the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,3,6,7,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,32,33,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,58,59,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $current = [string]$cell.Value2

    $parts = @($current.Split(",") | ForEach-Object { $_.Trim() })

    if ($parts.Count -gt 1) {
        $rotated = @($parts[1..($parts.Count - 1)]) + @($parts[0])
        $newValue = [string]::Join(", ", $rotated)
        $cell.Value2 = $newValue
    }
}
